$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 40.9
$ws.Range("I8").Value = 40.9
$ws.Range("K8").Value = 122.7
$ws.Range("M8").Value = 16.30000000000001

$ws.Range("H80").Value = 7860519
$ws.Range("I80").Value = 575.5625
$ws.Range("K80").Value = 1726.6875
$ws.Range("M80").Value = -728.6875

$ws.Range("H81").Value = 32328
$ws.Range("J81").Value = 32328
$ws.Range("L81").Value = 32328
$ws.Range("N81").Value = -34324

$ws.Range("H83").Value = 7860519
$ws.Range("I83").Value = 575.5625
$ws.Range("K83").Value = 5180.0625
$ws.Range("M83").Value = -188.0625

$ws.Range("H84").Value = 32328
$ws.Range("J84").Value = 32328
$ws.Range("L84").Value = 96984
$ws.Range("N84").Value = -106968

$ws.Range("H98").Value = 709.5185
$ws.Range("I98").Value = 762.8823
$ws.Range("J98").Value = 618.8
$ws.Range("K98").Value = 762.8823
$ws.Range("L98").Value = 618.8
$ws.Range("M98").Value = 735.1177
$ws.Range("N98").Value = -3614.8

$ws.Range("H122").Value = 709.5185
$ws.Range("I122").Value = 762.8823
$ws.Range("J122").Value = 618.8
$ws.Range("K122").Value = 2288.6469
$ws.Range("L122").Value = 1856.4
$ws.Range("M122").Value = 161.3531000000003
$ws.Range("N122").Value = -6756.4

$ws.Range("H137").Value = 27661.105
$ws.Range("I137").Value = 1303.2593
$ws.Range("J137").Value = 92357.63
$ws.Range("K137").Value = 3909.7779
$ws.Range("L137").Value = 277072.89
$ws.Range("M137").Value = -1359.7779
$ws.Range("N137").Value = -282172.89

$ws.Range("H138").Value = 1909.4884
$ws.Range("I138").Value = 866
$ws.Range("J138").Value = 2205.403
$ws.Range("K138").Value = 2598
$ws.Range("L138").Value = 6616.208999999999
$ws.Range("M138").Value = 2542
$ws.Range("N138").Value = -16896.209

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23280.037
$ws.Range("I32").Value = 27265.756
$ws.Range("K32").Value = 27265.756
$ws.Range("M32").Value = -26978.756

$ws.Range("H74").Value = 30303926
$ws.Range("I74").Value = 43478740
$ws.Range("J74").Value = 1850.5
$ws.Range("K74").Value = 43478740
$ws.Range("L74").Value = 1850.5
$ws.Range("M74").Value = -43477866
$ws.Range("N74").Value = -3598.5

$ws.Range("H77").Value = 30303926
$ws.Range("I77").Value = 43478740
$ws.Range("J77").Value = 1850.5
$ws.Range("K77").Value = 217393700
$ws.Range("L77").Value = 9252.5
$ws.Range("M77").Value = -217389332
$ws.Range("N77").Value = -17988.5

$ws.Range("H102").Value = 1338
$ws.Range("I102").Value = 1063.3334
$ws.Range("K102").Value = 1063.3334
$ws.Range("M102").Value = 558.6666

$ws.Range("H122").Value = 2950.6924
$ws.Range("I122").Value = 1845.5555
$ws.Range("J122").Value = 5437.25
$ws.Range("K122").Value = 5536.666499999999
$ws.Range("L122").Value = 16311.75
$ws.Range("M122").Value = -3086.666499999999
$ws.Range("N122").Value = -21211.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 42890
$ws.Range("J59").Value = 42890
$ws.Range("L59").Value = 42890
$ws.Range("N59").Value = -44584

$ws.Range("H86").Value = 1691.6061
$ws.Range("I86").Value = 1489.1154
$ws.Range("J86").Value = 2443.7144
$ws.Range("K86").Value = 1489.1154
$ws.Range("L86").Value = 2443.7144
$ws.Range("M86").Value = -366.1153999999999
$ws.Range("N86").Value = -4689.7144

$ws.Range("H89").Value = 1691.6061
$ws.Range("I89").Value = 1489.1154
$ws.Range("J89").Value = 2443.7144
$ws.Range("K89").Value = 7445.576999999999
$ws.Range("L89").Value = 12218.572
$ws.Range("M89").Value = -1829.576999999999
$ws.Range("N89").Value = -23450.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10010.293
$ws.Range("I31").Value = 18987.389
$ws.Range("K31").Value = 18987.389
$ws.Range("M31").Value = -18692.389

$ws.Range("H34").Value = 10010.293
$ws.Range("I34").Value = 18987.389
$ws.Range("K34").Value = 18987.389
$ws.Range("M34").Value = -18785.389

$ws.Range("H99").Value = 16132676
$ws.Range("I99").Value = 3520.7058
$ws.Range("J99").Value = 35718080
$ws.Range("K99").Value = 3520.7058
$ws.Range("L99").Value = 35718080
$ws.Range("M99").Value = -2022.7058
$ws.Range("N99").Value = -35721076

$ws.Range("H126").Value = 16132676
$ws.Range("I126").Value = 3520.7058
$ws.Range("J126").Value = 35718080
$ws.Range("K126").Value = 10562.1174
$ws.Range("L126").Value = 107154240
$ws.Range("M126").Value = -8092.117400000001
$ws.Range("N126").Value = -107159180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 232.5
$ws.Range("I17").Value = 158
$ws.Range("J17").Value = 285.7143
$ws.Range("K17").Value = 474
$ws.Range("L17").Value = 857.1428999999999
$ws.Range("M17").Value = -305
$ws.Range("N17").Value = -1195.1429

$ws.Range("H39").Value = 5666.3335
$ws.Range("J39").Value = 5666.3335
$ws.Range("L39").Value = 16999.0005
$ws.Range("N39").Value = -17587.0005

$ws.Range("H131").Value = 797.3099999999999
$ws.Range("J131").Value = 814.2371000000001
$ws.Range("L131").Value = 2442.7113
$ws.Range("N131").Value = -12522.7113

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 6171765
$ws.Range("J12").Value = 3010000
$ws.Range("L12").Value = 3010000
$ws.Range("N12").Value = -3010280

$ws.Range("H102").Value = 31251368
$ws.Range("I102").Value = 38462884
$ws.Range("K102").Value = 38462884
$ws.Range("M102").Value = -38461262

$ws.Range("H126").Value = 4836
$ws.Range("I126").Value = 3800
$ws.Range("J126").Value = 7037.5
$ws.Range("K126").Value = 11400
$ws.Range("L126").Value = 21112.5
$ws.Range("M126").Value = -8930
$ws.Range("N126").Value = -26052.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 4000
$ws.Range("J4").Value = 4000
$ws.Range("L4").Value = 4000
$ws.Range("N4").Value = -4226

$ws.Range("H28").Value = 4000
$ws.Range("J28").Value = 4000
$ws.Range("L28").Value = 4000
$ws.Range("N28").Value = -4464

$ws.Range("H37").Value = 4000
$ws.Range("J37").Value = 4000
$ws.Range("L37").Value = 4000
$ws.Range("N37").Value = -4214

$ws.Range("H40").Value = 7444.375
$ws.Range("I40").Value = 5333.3335
$ws.Range("J40").Value = 8711
$ws.Range("K40").Value = 5333.3335
$ws.Range("L40").Value = 8711
$ws.Range("M40").Value = -5197.3335
$ws.Range("N40").Value = -8983

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 52750
$ws.Range("J18").Value = 69333.336
$ws.Range("L18").Value = 69333.336
$ws.Range("N18").Value = -69679.336

$ws.Range("H122").Value = 2279.4
$ws.Range("I122").Value = 2099.5
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 6298.5
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -3848.5
$ws.Range("N122").Value = -13897
